$d = $word.ActiveDocument

# --- Remove the "_GoBack" bookmark that trails the final run of the
#     "D scores" bullet paragraph. ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- Insert a new bulleted "Highlights" list item after the paragraph
#     that ends in "... D scores" (it inherits the ListParagraph style
#     and numPr / rPr formatting of that paragraph automatically). ---
$srcPara = $d.Paragraphs.Item(4)
$srcPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item(5)
$newRange = $newPara.Range
$newRange.InsertAfter("Using an alternative scoring algorithm, the PI, did not improve performance")

$paraStart = $newPara.Range.Start

# Segment boundaries (start, end) within the new paragraph's text that
# must end up as separate <w:r> runs, matching the authored edit.
$segments = @(
    @(0, 10),
    @(10, 38),
    @(38, 39),
    @(39, 46),
    @(46, 47),
    @(47, 75)
)

# Force run breaks at each segment boundary by toggling a character
# property on/off over each segment (all segments share identical
# formatting, so the toggle leaves no visible / persisted trace, but
# it does split the underlying run at the segment boundaries).
foreach ($seg in $segments) {
    $segRange = $d.Range($paraStart + $seg[0], $paraStart + $seg[1])
    $segRange.Bold = 1
    $segRange.Bold = 0
}
